$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -21.945
    "A14"  = -21.659
    "A16"  = -21.985
    "A21"  = -20.022
    "A23"  = -20.302
    "A25"  = -21.775
    "A26"  = -21.319
    "A29"  = -21.236
    "A40"  = -20
    "A53"  = -21.948
    "A57"  = -22.213
    "A59"  = -22.5
    "A65"  = -21.533
    "A69"  = -21.602
    "A79"  = -21.167
    "A83"  = -22.035
    "A91"  = -21.533
    "A93"  = -21.22
    "A100" = -22.024
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
